$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.777.06"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.871.51"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'300.97"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5340"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("D8").Value = "'0.3739"
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "'0.07181"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'21.58"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.8894"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "'0.08187"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.885.95"
$ws.Range("E13").Value = "  +30.15%  "
$ws.Range("D14").Value = "'92.73"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "'5.306"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'14.82"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "'0.000008499"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "26.812.83"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "'4.981"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "'10.62"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'6.369"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'2.300"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "'146.08"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.06"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.729"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D29").Value = "'4.710"
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("D30").Value = "'4.614"
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").Value = "'0.09149"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'0.8044"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").Value = "'0.05013"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "'1.177"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'0.6113"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("D38").Value = "'3.194"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("D39").Value = "'0.01951"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").Value = "'1.066"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'6.573"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'0.5267"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").Value = "'8.777"
$ws.Range("D44").Value = "'114.61"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'0.1492"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'1.656"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "'9.889"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "'37.54"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").Value = "'0.06058"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'62.01"
$ws.Range("E51").Value = "  -3.70%  "
